$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Operator2 value (H2): "Equal" -> "Does Not Contains"
$ws.Range("H2").Value = "Does Not Contains"

# FilterValue1 (F2): 233 -> 23
$ws.Range("F2").Value = 23

# Column H widens to fit the new, longer operator text
$ws.Columns("H").ColumnWidth = 15.5

# Active cell returns to the top-left (default) after the edits
$ws.Range("A1").Select() | Out-Null
